# 7.8 History Card & Advanced Story
# Rewrites Fu's interrogation answers (rows 4-19) with the revised dialogue
# text from the commit, keeps row 20 ("Goto" / "StoryScript9") intact,
# updates the two rows whose wrapped text got longer/shorter (row heights),
# and moves the sheet selection to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = "I used to live in JiuJiang county at the foot of the mountain."
$ws.Range("B5").Value  = "Alas, my family was impoverished, and I possessed no particular skills——only a talent for poetry and calligraphy."
$ws.Range("B6").Value  = "A month ago, I happened to meet Ming in downtown. We got along very well, and he invited me to reside at the manor."
$ws.Range("B7").Value  = "Since then, I’ve occasionally taught him a bit of poetry."
$ws.Range("B8").Value  = "When was the last time you saw the Lord?"
$ws.Range("B9").Value  = "It may have been several days ago."
$ws.Range("B10").Value = "Why is that? Don’t you often see the Lord, living here at the manor?"
$ws.Range("B11").Value = "As you may know, the Lord hoped for Ming to inherit an official position."
$ws.Range("B12").Value = "He disapproved of Ming studying poetry and the arts, and thus did not look kindly upon me."
$ws.Range("B13").Value = "Please describe your actions around the time of the evening banquet."
$ws.Range("B14").Value = "To avoid any unpleasantness, I remained alone in my room the whole time. I didn’t attend the banquet."
$ws.Range("B15").Value = "Did you interact with anyone else today?"
$ws.Range("B16").Value = "After around 3 PM, when most had likely finished their lunch, I went to the canteen."
$ws.Range("B17").Value = "On the way back, I happened to see you and Ming chatting in the garden."
$ws.Range("B18").Value = "Around 7.45 PM, Ming suddenly knocked on my door, asking if I knew where the Lord was."
$ws.Range("B19").Value = "I answered truthfully, and he left shortly after."

# Row 20 keeps the same text ("Goto" / "StoryScript9") - only its shared
# string index shifts because of the sst reshuffle above, so no visible
# change, but we re-assert it for good measure.
$ws.Range("A20").Value = "Goto"
$ws.Range("B20").Value = "StoryScript9"

# Row heights: row 4 and row 18 text got shorter (wrap now needs fewer
# lines) while row 5 got longer (needs one more wrapped line).
$ws.Rows.Item(4).RowHeight  = 34
$ws.Rows.Item(5).RowHeight  = 51
$ws.Rows.Item(18).RowHeight = 34

# Move the saved selection/scroll position to match the new view.
$ws.Range("C22").Select()
